$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on F2 first so Rows.Insert() below doesn't
# leave a stale hyperlink object pointing at the (about to be emptied) F2 cell.
if ($ws.Range("F2").Hyperlinks.Count -gt 0) {
    $ws.Range("F2").Hyperlinks.Delete()
}

# Insert a brand-new row above the current row 2 - this pushes the existing
# "WordPress" job row (and its formatting, incl. the Hyperlink cell style)
# down to row 3.
$ws.Rows(2).Insert()

# --- Row 2: new job entry fetched at 2025-11-23 12:31:24 ---
$ws.Range("A2").Value = "2025-11-23 12:31:24"
$ws.Range("B2").Value = "【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5439921"
$ws.Range("G2").Value = 155
$ws.Range("H2").Value = "★Java ◆開発"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5439921") | Out-Null
$ws.Range("F2").Style = "Hyperlink"

# --- Row 3: pre-existing job entry, re-stamped with the same refresh time ---
$ws.Range("A3").Value = "2025-11-23 12:31:24"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439670") | Out-Null
$ws.Range("F3").Style = "Hyperlink"

# Widen column D to fit the new, longer price string.
# (the engine stores <col width> as ColumnWidth + 5/6, so back that offset out
# here to land on an exact stored width of 30)
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
